# Update "想去人数" (column F) counts across sheets to the newly scraped
# values, as published to gh-pages at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# -- 展览 (Exhibitions) -------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 6961
$wsExpo.Range("F4").Value = 56
$wsExpo.Range("F6").Value = 157
$wsExpo.Range("F7").Value = 6824
$wsExpo.Range("F8").Value = 73
$wsExpo.Range("F9").Value = 200
$wsExpo.Range("F10").Value = 1287
$wsExpo.Range("F12").Value = 0
$wsExpo.Range("F13").Value = 406
$wsExpo.Range("F15").Value = 17
$wsExpo.Range("F17").Value = 48
$wsExpo.Range("F19").Value = 15
$wsExpo.Range("F20").Value = 5205
$wsExpo.Range("F21").Value = 117
$wsExpo.Range("F22").Value = 161
$wsExpo.Range("F23").Value = 626
$wsExpo.Range("F25").Value = 229

# -- 演出 (Performances) -------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 46

# -- 全部类型 (All types) -------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6961
$wsAll.Range("F6").Value = 157
$wsAll.Range("F7").Value = 6824
$wsAll.Range("F9").Value = 200
$wsAll.Range("F10").Value = 1287
$wsAll.Range("F11").Value = 0
$wsAll.Range("F12").Value = 108
$wsAll.Range("F13").Value = 406
$wsAll.Range("F17").Value = 0
$wsAll.Range("F21").Value = 5205
$wsAll.Range("F22").Value = 46
$wsAll.Range("F27").Value = 229
